$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 939.6923
$ws.Range("I2").Value = 252.6
$ws.Range("J2").Value = 3230
$ws.Range("K2").Value = 252.6
$ws.Range("L2").Value = 3230
$ws.Range("M2").Value = -139.6
$ws.Range("N2").Value = -3456
$ws.Range("H28").Value = 1158
$ws.Range("I28").Value = 1305.8182
$ws.Range("K28").Value = 1305.8182
$ws.Range("M28").Value = -820.8181999999999
$ws.Range("H40").Value = 2075.8823
$ws.Range("I40").Value = 1644
$ws.Range("J40").Value = 2561.75
$ws.Range("K40").Value = 1644
$ws.Range("L40").Value = 2561.75
$ws.Range("M40").Value = -1469
$ws.Range("N40").Value = -2911.75
$ws.Range("H111").Value = 797
$ws.Range("I111").Value = 796
$ws.Range("J111").Value = 798
$ws.Range("K111").Value = 2388
$ws.Range("L111").Value = 2394
$ws.Range("M111").Value = 679
$ws.Range("N111").Value = -8528
$ws.Range("H116").Value = 6999.5
$ws.Range("I116").Value = 6666
$ws.Range("K116").Value = 6666
$ws.Range("M116").Value = -3224
$ws.Range("H135").Value = 392.46155
$ws.Range("I135").Value = 450.9
$ws.Range("J135").Value = 197.66667
$ws.Range("K135").Value = 4058.1
$ws.Range("L135").Value = 1779.00003
$ws.Range("M135").Value = -1523.1
$ws.Range("N135").Value = -6849.00003
$ws.Range("H141").Value = 2935.5789
$ws.Range("I141").Value = 2866.5334
$ws.Range("K141").Value = 8599.600199999999
$ws.Range("M141").Value = -3419.600199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3083338.2
$ws.Range("I32").Value = 4670900
$ws.Range("J32").Value = 701995.8
$ws.Range("K32").Value = 4670900
$ws.Range("L32").Value = 701995.8
$ws.Range("M32").Value = -4670613
$ws.Range("N32").Value = -702569.8
$ws.Range("H45").Value = 5282.3335
$ws.Range("I45").Value = 6035.8
$ws.Range("K45").Value = 6035.8
$ws.Range("M45").Value = -5658.8
$ws.Range("H88").Value = 1739.1428
$ws.Range("I88").Value = 1634.25
$ws.Range("J88").Value = 1879
$ws.Range("K88").Value = 1634.25
$ws.Range("L88").Value = 1879
$ws.Range("M88").Value = -1228.25
$ws.Range("N88").Value = -2691
$ws.Range("H91").Value = 1739.1428
$ws.Range("I91").Value = 1634.25
$ws.Range("J91").Value = 1879
$ws.Range("K91").Value = 1634.25
$ws.Range("L91").Value = 1879
$ws.Range("M91").Value = -230.25
$ws.Range("N91").Value = -4687
$ws.Range("H113").Value = 129398
$ws.Range("J113").Value = 129398
$ws.Range("L113").Value = 129398
$ws.Range("N113").Value = -138076
$ws.Range("H132").Value = 2735.3333
$ws.Range("I132").Value = 2735.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8205.999899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5675.999899999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1608.1818
$ws.Range("I86").Value = 1554.4445
$ws.Range("J86").Value = 1850
$ws.Range("K86").Value = 1554.4445
$ws.Range("L86").Value = 1850
$ws.Range("M86").Value = -431.4445000000001
$ws.Range("N86").Value = -4096
$ws.Range("H89").Value = 1608.1818
$ws.Range("I89").Value = 1554.4445
$ws.Range("J89").Value = 1850
$ws.Range("K89").Value = 7772.2225
$ws.Range("L89").Value = 9250
$ws.Range("M89").Value = -2156.2225
$ws.Range("N89").Value = -20482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 8000
$ws.Range("I93").Value = 8000
$ws.Range("K93").Value = 8000
$ws.Range("M93").Value = -6128
$ws.Range("H103").Value = 6524
$ws.Range("I103").Value = 6524
$ws.Range("K103").Value = 6524
$ws.Range("M103").Value = -5352
$ws.Range("H105").Value = 2502.9167
$ws.Range("I105").Value = 1899.8
$ws.Range("J105").Value = 2933.7144
$ws.Range("K105").Value = 1899.8
$ws.Range("L105").Value = 2933.7144
$ws.Range("M105").Value = -152.8
$ws.Range("N105").Value = -6427.7144
$ws.Range("H107").Value = 367.5
$ws.Range("I107").Value = 338.66666
$ws.Range("J107").Value = 396.33334
$ws.Range("K107").Value = 338.66666
$ws.Range("L107").Value = 396.33334
$ws.Range("M107").Value = 1581.33334
$ws.Range("N107").Value = -4236.33334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 148.09091
$ws.Range("J2").Value = 200.42857
$ws.Range("L2").Value = 1202.57142
$ws.Range("N2").Value = -1428.57142
$ws.Range("H32").Value = 1700
$ws.Range("I32").Value = 900
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 2700
$ws.Range("L32").Value = 7500
$ws.Range("M32").Value = -2417
$ws.Range("N32").Value = -8066
$ws.Range("H37").Value = 69979.664
$ws.Range("J37").Value = 69979.664
$ws.Range("L37").Value = 209938.992
$ws.Range("N37").Value = -210162.992
$ws.Range("H40").Value = 71
$ws.Range("J40").Value = 82.666664
$ws.Range("L40").Value = 330.666656
$ws.Range("N40").Value = -468.666656
$ws.Range("H121").Value = 8823.368
$ws.Range("I121").Value = 22144
$ws.Range("K121").Value = 66432
$ws.Range("M121").Value = -65122
$ws.Range("H130").Value = 1699
$ws.Range("I130").Value = 1699
$ws.Range("K130").Value = 5097
$ws.Range("M130").Value = -77
$ws.Range("H131").Value = 418918.25
$ws.Range("I131").Value = 1118.8334
$ws.Range("K131").Value = 3356.5002
$ws.Range("M131").Value = 1683.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 112950
$ws.Range("J5").Value = 900
$ws.Range("L5").Value = 900
$ws.Range("N5").Value = -1124
$ws.Range("H70").Value = 6349.75
$ws.Range("I70").Value = 6333
$ws.Range("J70").Value = 6400
$ws.Range("K70").Value = 6333
$ws.Range("L70").Value = 6400
$ws.Range("M70").Value = -6063
$ws.Range("N70").Value = -6940
$ws.Range("H73").Value = 6349.75
$ws.Range("I73").Value = 6333
$ws.Range("J73").Value = 6400
$ws.Range("K73").Value = 6333
$ws.Range("L73").Value = 6400
$ws.Range("M73").Value = -5397
$ws.Range("N73").Value = -8272
$ws.Range("H113").Value = 239.4
$ws.Range("I113").Value = 239.4
$ws.Range("K113").Value = 239.4
$ws.Range("M113").Value = 1930.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 500
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 28999.334
$ws.Range("J74").Value = 28999.334
$ws.Range("L74").Value = 28999.334
$ws.Range("N74").Value = -30871.334
$ws.Range("H77").Value = 28999.334
$ws.Range("J77").Value = 28999.334
$ws.Range("L77").Value = 86998.00199999999
$ws.Range("N77").Value = -96358.00199999999
$ws.Range("H136").Value = 2960.739
$ws.Range("I136").Value = 3090.4285
$ws.Range("K136").Value = 9271.2855
$ws.Range("M136").Value = -6721.2855
